$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column C filled with 1 for rows 2-5 (SL_No column)
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 1
$ws.Range("C4").Value = 1
$ws.Range("C5").Value = 1

# Rename header labels to uppercase canonical names
$ws.Range("U1").Value = "INGREDIENT_NAME"
$ws.Range("V1").Value = "WEIGHT_IN_G"

# Set column widths to match bestFit sizing used after the edit
$ws.Columns("U").ColumnWidth = 17.711495535714285
$ws.Columns("V").ColumnWidth = 12.711495535714286

# Update selection to reflect where the user left off
[void]$ws.Range("C6").Select()
